$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-43) holds the "Förändrad" date serial value.
# Update it from 45812 (2025-06-04) to 45813 (2025-06-05) for every data row.
$range = $ws.Range("C2:C43")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45812) {
        $cell.Value2 = 45813
    }
}
